# Edit script applying the "river update May 2024" refresh to the trend-results sheet.
# The whole data table (site x parameter x trend-period) is refreshed with new
# statistics, and a new trend-period=15 block (rows 28-37) is appended.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update changed cells in existing rows 2-27 ---
# Row 2
$ws.Range("F2").Value2 = 0.096859914106757
$ws.Range("G2").Value2 = 0.0952380952380952
$ws.Range("H2").Value2 = 0.952380952380952
$ws.Range("J2").Value2 = 0.46
$ws.Range("K2").Value2 = -0.183904993815708
$ws.Range("L2").Value2 = -0.414705783389464
$ws.Range("M2").Value2 = 0.0335781061833756
$ws.Range("N2").Value2 = -39.9793464816757
$ws.Range("P2").Value2 = "Very unlikely improving"

# Row 3
$ws.Range("F3").Value2 = 0.017182123185329
$ws.Range("H3").Value2 = 0.895833333333333
$ws.Range("J3").Value2 = 10.815
$ws.Range("K3").Value2 = -0.193953242481207
$ws.Range("L3").Value2 = -0.370544215825512
$ws.Range("M3").Value2 = -0.046994455956346
$ws.Range("N3").Value2 = -1.79337256108374
$ws.Range("P3").Value2 = "Extremely unlikely increasing"

# Row 4
$ws.Range("F4").Value2 = 0.344770659432582
$ws.Range("G4").Value2 = 0.291666666666667
$ws.Range("H4").Value2 = 0.229166666666667
$ws.Range("J4").Value2 = 0.0065
$ws.Range("L4").Value2 = 0
$ws.Range("M4").Value2 = 0.0007164493322046
$ws.Range("P4").Value2 = "As likely as not improving"

# Row 5
$ws.Range("F5").Value2 = 0.895014086955298
$ws.Range("G5").Value2 = 0.0208333333333333
$ws.Range("H5").Value2 = 0.833333333333333
$ws.Range("J5").Value2 = 135
$ws.Range("K5").Value2 = -13.4071827862405
$ws.Range("L5").Value2 = -41.2445147976945
$ws.Range("M5").Value2 = 3.96347015156552
$ws.Range("N5").Value2 = -9.93124650832628
$ws.Range("P5").Value2 = "Likely improving"

# Row 6
$ws.Range("E6").Value2 = "WARNING: Sen slope based on two censored values"
$ws.Range("F6").Value2 = 0.273006084795537
$ws.Range("G6").Value2 = 0.644444444444444
$ws.Range("H6").Value2 = 0.488888888888889
$ws.Range("I6").Value2 = 7

# Row 7
$ws.Range("E7").Value2 = "ok"
$ws.Range("F7").Value2 = 0.999996792093399
$ws.Range("G7").Value2 = 0.0416666666666667
$ws.Range("H7").Value2 = 0.1875
$ws.Range("J7").Value2 = 0.0045
$ws.Range("K7").Value2 = -0.0009221564289432
$ws.Range("L7").Value2 = -0.0012252342713039
$ws.Range("M7").Value2 = -0.0006363240418118
$ws.Range("N7").Value2 = -20.492365087628
$ws.Range("P7").Value2 = "Virtually certain improving"

# Row 8
$ws.Range("F8").Value2 = 0.143074675100918
$ws.Range("H8").Value2 = 0.979166666666667
$ws.Range("J8").Value2 = 0.631
$ws.Range("K8").Value2 = 0.02435
$ws.Range("L8").Value2 = -0.0191760819886591
$ws.Range("M8").Value2 = 0.0561892874243412
$ws.Range("N8").Value2 = 3.85895404120444
$ws.Range("P8").Value2 = "Unlikely improving"

# Row 9
$ws.Range("E9").Value2 = "ok"
$ws.Range("F9").Value2 = 0.049345077881456
$ws.Range("H9").Value2 = 0.829787234042553
$ws.Range("J9").Value2 = 8.07
$ws.Range("K9").Value2 = -0.0565791738382099
$ws.Range("L9").Value2 = -0.11400531801377
$ws.Range("M9").Value2 = 0
$ws.Range("N9").Value2 = -0.701105004190953
$ws.Range("P9").Value2 = "Extremely unlikely increasing"

# Row 10
$ws.Range("F10").Value2 = 0.164105656945935
$ws.Range("H10").Value2 = 0.979166666666667
$ws.Range("J10").Value2 = 0.647
$ws.Range("K10").Value2 = 0.0243651789723073
$ws.Range("L10").Value2 = -0.0201370695082495
$ws.Range("M10").Value2 = 0.05705309955189
$ws.Range("N10").Value2 = 3.76587001117578
$ws.Range("P10").Value2 = "Unlikely improving"

# Row 11
$ws.Range("B11").Value2 = "Total Nitrogen"
$ws.Range("D11").Value2 = $false
$ws.Range("F11").Value2 = 0.367705366050189
$ws.Range("G11").Value2 = 0
$ws.Range("H11").Value2 = 0.770833333333333
$ws.Range("I11").Value2 = 0
$ws.Range("J11").Value2 = 0.8
$ws.Range("K11").Value2 = 0.0090807138533379
$ws.Range("L11").Value2 = -0.0321903461122343
$ws.Range("M11").Value2 = 0.0487973844492751
$ws.Range("N11").Value2 = 1.13508923166724
$ws.Range("P11").Value2 = "As likely as not improving"
$ws.Range("W11").Value2 = "g/m3"

# Row 12
$ws.Range("B12").Value2 = "Total Phosphorus"
$ws.Range("D12").Value2 = $false
$ws.Range("F12").Value2 = 0.72471266978495
$ws.Range("H12").Value2 = 0.583333333333333
$ws.Range("J12").Value2 = 0.021
$ws.Range("K12").Value2 = -0.0005006549536519
$ws.Range("L12").Value2 = -0.0025912275998017
$ws.Range("M12").Value2 = 0.0009761207165859
$ws.Range("N12").Value2 = -2.38407120786649
$ws.Range("P12").Value2 = "Likely improving"

# Row 13
$ws.Range("B13").Value2 = "Turbidity"
$ws.Range("F13").Value2 = 0.941133111119049
$ws.Range("H13").Value2 = 0.979166666666667
$ws.Range("J13").Value2 = 2.245
$ws.Range("K13").Value2 = -0.309285322720736
$ws.Range("L13").Value2 = -0.909171063966783
$ws.Range("M13").Value2 = 0.0237769653030699
$ws.Range("N13").Value2 = -13.7766290744203
$ws.Range("P13").Value2 = "Very likely improving"
$ws.Range("W13").Value2 = "NTU/FNU"

# Row 14
$ws.Range("B14").Value2 = "Visual Clarity"
$ws.Range("C14").Value2 = 10
$ws.Range("D14").Value2 = $false
$ws.Range("F14").Value2 = 0.010618377183117
$ws.Range("G14").Value2 = 0.028169014084507
$ws.Range("H14").Value2 = 0.816901408450704
$ws.Range("I14").Value2 = 2
$ws.Range("J14").Value2 = 0.9
$ws.Range("K14").Value2 = -0.100371279761905
$ws.Range("L14").Value2 = -0.178694570320577
$ws.Range("M14").Value2 = -0.0268298765481437
$ws.Range("N14").Value2 = -11.1523644179894
$ws.Range("P14").Value2 = "Extremely unlikely improving"
$ws.Range("W14").Value2 = "m"

# Row 15
$ws.Range("B15").Value2 = "Dissolved Oxygen Concentration"
$ws.Range("D15").Value2 = $true
$ws.Range("F15").Value2 = 0.776519453733732
$ws.Range("G15").Value2 = 0
$ws.Range("H15").Value2 = 0.830188679245283
$ws.Range("I15").Value2 = 0
$ws.Range("J15").Value2 = 10.72
$ws.Range("K15").Value2 = 0.0198775510204086
$ws.Range("L15").Value2 = -0.0375994753898634
$ws.Range("M15").Value2 = 0.0719073109103431
$ws.Range("N15").Value2 = 0.185424916235155
$ws.Range("P15").Value2 = "Likely increasing"
$ws.Range("W15").Value2 = "g/m3"

# Row 16
$ws.Range("B16").Value2 = "Dissolved Reactive Phosphorus"
$ws.Range("D16").Value2 = $false
$ws.Range("E16").Value2 = "WARNING: Sen slope influenced by censored values"
$ws.Range("F16").Value2 = 0.878094644072962
$ws.Range("G16").Value2 = 0.242990654205607
$ws.Range("H16").Value2 = 0.149532710280374
$ws.Range("I16").Value2 = 1
$ws.Range("J16").Value2 = 0.008
$ws.Range("K16").Value2 = 0
$ws.Range("L16").Value2 = -0.0003783827779132
$ws.Range("M16").Value2 = 0
$ws.Range("N16").Value2 = 0
$ws.Range("P16").Value2 = "Likely improving"
$ws.Range("W16").Value2 = "mg/L"

# Row 17
$ws.Range("B17").Value2 = "E. coli"
$ws.Range("E17").Value2 = "ok"
$ws.Range("F17").Value2 = 0.923972740061657
$ws.Range("G17").Value2 = 0.0093457943925233
$ws.Range("H17").Value2 = 0.747663551401869
$ws.Range("J17").Value2 = 160
$ws.Range("K17").Value2 = -6.90598739495798
$ws.Range("L17").Value2 = -14.7962802636473
$ws.Range("M17").Value2 = 1.00332593730069
$ws.Range("N17").Value2 = -4.31624212184874
$ws.Range("P17").Value2 = "Very likely improving"
$ws.Range("W17").Value2 = "E. coli/100 mL"

# Row 18
$ws.Range("B18").Value2 = "Ammoniacal Nitrogen (NH4)"
$ws.Range("E18").Value2 = "WARNING: Sen slope influenced by censored values"
$ws.Range("F18").Value2 = 0.0669648863040605
$ws.Range("G18").Value2 = 0.762376237623762
$ws.Range("H18").Value2 = 0.376237623762376
$ws.Range("I18").Value2 = 17
$ws.Range("J18").Value2 = 0.005
$ws.Range("K18").Value2 = 0
$ws.Range("L18").Value2 = 0
$ws.Range("M18").Value2 = 0
$ws.Range("N18").Value2 = 0
$ws.Range("P18").Value2 = "Very unlikely improving"
$ws.Range("W18").Value2 = "mg/L"

# Row 19
$ws.Range("B19").Value2 = "Nitrite Nitrogen (NO2)"
$ws.Range("E19").Value2 = "ok"
$ws.Range("F19").Value2 = 0.999723028200602
$ws.Range("G19").Value2 = 0.0186915887850467
$ws.Range("H19").Value2 = 0.149532710280374
$ws.Range("I19").Value2 = 1
$ws.Range("K19").Value2 = -0.0002595948827292
$ws.Range("L19").Value2 = -0.0004141156462585
$ws.Range("M19").Value2 = -0.0001119792742978
$ws.Range("N19").Value2 = -5.19189765458422
$ws.Range("P19").Value2 = "Virtually certain improving"

# Row 20
$ws.Range("B20").Value2 = "Nitrate Nitrogen (NO3)"
$ws.Range("D20").Value2 = $true
$ws.Range("E20").Value2 = "WARNING: Sen slope based on tied non-censored values"
$ws.Range("F20").Value2 = 0.5
$ws.Range("G20").Value2 = 0
$ws.Range("H20").Value2 = 0.97196261682243
$ws.Range("I20").Value2 = 0
$ws.Range("J20").Value2 = 0.66
$ws.Range("K20").Value2 = 0
$ws.Range("L20").Value2 = -0.0157264999912368
$ws.Range("M20").Value2 = 0.0073672729693228
$ws.Range("N20").Value2 = 0
$ws.Range("P20").Value2 = "As likely as not improving"

# Row 21
$ws.Range("B21").Value2 = "pH"
$ws.Range("D21").Value2 = $false
$ws.Range("F21").Value2 = 0.907673882045815
$ws.Range("H21").Value2 = 0.685714285714286
$ws.Range("J21").Value2 = 7.92
$ws.Range("K21").Value2 = 0.0191355681280149
$ws.Range("L21").Value2 = -0.0036235119047618
$ws.Range("M21").Value2 = 0.0391566852992791
$ws.Range("N21").Value2 = 0.241610708687057
$ws.Range("P21").Value2 = "Very likely increasing"
$ws.Range("W21").Value2 = ""

# Row 22
$ws.Range("B22").Value2 = "SIN (Soluble Inorganic nitrogen)"
$ws.Range("D22").Value2 = $true
$ws.Range("F22").Value2 = 0.440281460095996
$ws.Range("H22").Value2 = 0.878504672897196
$ws.Range("J22").Value2 = 0.676
$ws.Range("K22").Value2 = 0.0005702576112412
$ws.Range("L22").Value2 = -0.0151839697940833
$ws.Range("M22").Value2 = 0.0079952348571241
$ws.Range("N22").Value2 = 0.0843576347989956
$ws.Range("P22").Value2 = "As likely as not improving"
$ws.Range("W22").Value2 = "g/m3"

# Row 23
$ws.Range("B23").Value2 = "Total Nitrogen"
$ws.Range("F23").Value2 = 0.873798018537367
$ws.Range("H23").Value2 = 0.663551401869159
$ws.Range("J23").Value2 = 0.86
$ws.Range("K23").Value2 = -0.0061680739388407
$ws.Range("L23").Value2 = -0.0165097687234024
$ws.Range("M23").Value2 = 0.0018885149786085
$ws.Range("N23").Value2 = -0.7172178998652
$ws.Range("P23").Value2 = "Likely improving"

# Row 24
$ws.Range("B24").Value2 = "Total Phosphorus"
$ws.Range("E24").Value2 = "ok"
$ws.Range("F24").Value2 = 0.926979124249585
$ws.Range("G24").Value2 = 0
$ws.Range("H24").Value2 = 0.429906542056075
$ws.Range("I24").Value2 = 0
$ws.Range("J24").Value2 = 0.024
$ws.Range("K24").Value2 = -0.000501717032967
$ws.Range("L24").Value2 = -0.001003434065934
$ws.Range("N24").Value2 = -2.09048763736263
$ws.Range("P24").Value2 = "Very likely improving"
$ws.Range("W24").Value2 = "g/m3"

# Row 25
$ws.Range("B25").Value2 = "Turbidity"
$ws.Range("E25").Value2 = "ok"
$ws.Range("F25").Value2 = 0.808254698496209
$ws.Range("H25").Value2 = 0.934579439252336
$ws.Range("J25").Value2 = 3.97
$ws.Range("K25").Value2 = -0.0643297455968689
$ws.Range("L25").Value2 = -0.222772362141085
$ws.Range("M25").Value2 = 0.0542680267301122
$ws.Range("N25").Value2 = -1.6203966145307
$ws.Range("P25").Value2 = "Likely improving"
$ws.Range("W25").Value2 = "NTU/FNU"

# Row 26
$ws.Range("B26").Value2 = "Visual Clarity"
$ws.Range("C26").Value2 = 15
$ws.Range("F26").Value2 = 0.151761419627249
$ws.Range("G26").Value2 = 0.0192307692307692
$ws.Range("H26").Value2 = 0.663461538461538
$ws.Range("I26").Value2 = 2
$ws.Range("J26").Value2 = 0.91
$ws.Range("K26").Value2 = -0.0260892857142857
$ws.Range("L26").Value2 = -0.0733289382624688
$ws.Range("M26").Value2 = 0.0099760819089
$ws.Range("N26").Value2 = -2.86695447409733
$ws.Range("P26").Value2 = "Unlikely improving"
$ws.Range("W26").Value2 = "m"

# Row 27
$ws.Range("B27").Value2 = "Dissolved Oxygen Concentration"
$ws.Range("C27").Value2 = 15
$ws.Range("F27").Value2 = 0.715590382395065
$ws.Range("H27").Value2 = 0.806666666666667
$ws.Range("J27").Value2 = 10.72
$ws.Range("K27").Value2 = 0.0099795081967213
$ws.Range("L27").Value2 = -0.0284748302125065
$ws.Range("M27").Value2 = 0.0401373626373625
$ws.Range("N27").Value2 = 0.0930924272082213
$ws.Range("P27").Value2 = "Likely increasing"
$ws.Range("W27").Value2 = "g/m3"

# --- Append new rows 28-37 (trend period = 15) ---
# Row 28
$ws.Range("A28").Value2 = "Tiraumea u/s Manawatu Confluence"
$ws.Range("B28").Value2 = "Dissolved Reactive Phosphorus"
$ws.Range("C28").Value2 = 15
$ws.Range("D28").Value2 = $false
$ws.Range("E28").Value2 = "WARNING: Sen slope influenced by censored values"
$ws.Range("F28").Value2 = 0.944754591566321
$ws.Range("G28").Value2 = 0.225165562913907
$ws.Range("H28").Value2 = 0.132450331125828
$ws.Range("I28").Value2 = 1
$ws.Range("J28").Value2 = 0.008
$ws.Range("K28").Value2 = 0
$ws.Range("L28").Value2 = -0.0002933266044311
$ws.Range("M28").Value2 = 0
$ws.Range("N28").Value2 = 0
$ws.Range("O28").Value2 = "RepSite"
$ws.Range("P28").Value2 = "Very likely improving"
$ws.Range("Q28").Value2 = 1845196
$ws.Range("R28").Value2 = 5525095
$ws.Range("S28").Value2 = "Tararua District"
$ws.Range("T28").Value2 = "Manawatū"
$ws.Range("U28").Value2 = "Tiraumea"
$ws.Range("V28").Value2 = "Mana_7b"
$ws.Range("W28").Value2 = "mg/L"

# Row 29
$ws.Range("A29").Value2 = "Tiraumea u/s Manawatu Confluence"
$ws.Range("B29").Value2 = "E. coli"
$ws.Range("C29").Value2 = 15
$ws.Range("D29").Value2 = $false
$ws.Range("E29").Value2 = "ok"
$ws.Range("F29").Value2 = 0.969681970328645
$ws.Range("G29").Value2 = 0.0066225165562913
$ws.Range("H29").Value2 = 0.748344370860927
$ws.Range("I29").Value2 = 1
$ws.Range("J29").Value2 = 170
$ws.Range("K29").Value2 = -6.07027934936351
$ws.Range("L29").Value2 = -11.4100224586356
$ws.Range("M29").Value2 = -0.706917139684145
$ws.Range("N29").Value2 = -3.57075255844912
$ws.Range("O29").Value2 = "RepSite"
$ws.Range("P29").Value2 = "Extremely likely improving"
$ws.Range("Q29").Value2 = 1845196
$ws.Range("R29").Value2 = 5525095
$ws.Range("S29").Value2 = "Tararua District"
$ws.Range("T29").Value2 = "Manawatū"
$ws.Range("U29").Value2 = "Tiraumea"
$ws.Range("V29").Value2 = "Mana_7b"
$ws.Range("W29").Value2 = "E. coli/100 mL"

# Row 30
$ws.Range("A30").Value2 = "Tiraumea u/s Manawatu Confluence"
$ws.Range("B30").Value2 = "Ammoniacal Nitrogen (NH4)"
$ws.Range("C30").Value2 = 15
$ws.Range("D30").Value2 = $false
$ws.Range("E30").Value2 = "WARNING: Sen slope influenced by censored values"
$ws.Range("F30").Value2 = 0.863711141749014
$ws.Range("G30").Value2 = 0.703448275862069
$ws.Range("H30").Value2 = 0.4
$ws.Range("I30").Value2 = 18
$ws.Range("J30").Value2 = 0.005
$ws.Range("K30").Value2 = 0
$ws.Range("L30").Value2 = 0
$ws.Range("M30").Value2 = 0
$ws.Range("N30").Value2 = 0
$ws.Range("O30").Value2 = "RepSite"
$ws.Range("P30").Value2 = "Likely improving"
$ws.Range("Q30").Value2 = 1845196
$ws.Range("R30").Value2 = 5525095
$ws.Range("S30").Value2 = "Tararua District"
$ws.Range("T30").Value2 = "Manawatū"
$ws.Range("U30").Value2 = "Tiraumea"
$ws.Range("V30").Value2 = "Mana_7b"
$ws.Range("W30").Value2 = "mg/L"

# Row 31
$ws.Range("A31").Value2 = "Tiraumea u/s Manawatu Confluence"
$ws.Range("B31").Value2 = "Nitrite Nitrogen (NO2)"
$ws.Range("C31").Value2 = 15
$ws.Range("D31").Value2 = $false
$ws.Range("E31").Value2 = "WARNING: Sen slope influenced by censored values"
$ws.Range("F31").Value2 = 0.377289937424923
$ws.Range("G31").Value2 = 0.152317880794702
$ws.Range("H31").Value2 = 0.105960264900662
$ws.Range("I31").Value2 = 3
$ws.Range("J31").Value2 = 0.004
$ws.Range("K31").Value2 = 0
$ws.Range("L31").Value2 = 0
$ws.Range("M31").Value2 = 0.0000921561175295943
$ws.Range("N31").Value2 = 0
$ws.Range("O31").Value2 = "RepSite"
$ws.Range("P31").Value2 = "As likely as not improving"
$ws.Range("Q31").Value2 = 1845196
$ws.Range("R31").Value2 = 5525095
$ws.Range("S31").Value2 = "Tararua District"
$ws.Range("T31").Value2 = "Manawatū"
$ws.Range("U31").Value2 = "Tiraumea"
$ws.Range("V31").Value2 = "Mana_7b"
$ws.Range("W31").Value2 = "mg/L"

# Row 32
$ws.Range("A32").Value2 = "Tiraumea u/s Manawatu Confluence"
$ws.Range("B32").Value2 = "Nitrate Nitrogen (NO3)"
$ws.Range("C32").Value2 = 15
$ws.Range("D32").Value2 = $true
$ws.Range("E32").Value2 = "ok"
$ws.Range("F32").Value2 = 0.529029969330229
$ws.Range("G32").Value2 = 0
$ws.Range("H32").Value2 = 0.927152317880795
$ws.Range("I32").Value2 = 0
$ws.Range("J32").Value2 = 0.6524
$ws.Range("K32").Value2 = -0.0002496582365003
$ws.Range("L32").Value2 = -0.0079510341495895
$ws.Range("M32").Value2 = 0.0052451162226004
$ws.Range("N32").Value2 = -0.0382676634733817
$ws.Range("O32").Value2 = "RepSite"
$ws.Range("P32").Value2 = "As likely as not improving"
$ws.Range("Q32").Value2 = 1845196
$ws.Range("R32").Value2 = 5525095
$ws.Range("S32").Value2 = "Tararua District"
$ws.Range("T32").Value2 = "Manawatū"
$ws.Range("U32").Value2 = "Tiraumea"
$ws.Range("V32").Value2 = "Mana_7b"
$ws.Range("W32").Value2 = "mg/L"

# Row 33
$ws.Range("A33").Value2 = "Tiraumea u/s Manawatu Confluence"
$ws.Range("B33").Value2 = "pH"
$ws.Range("C33").Value2 = 15
$ws.Range("D33").Value2 = $false
$ws.Range("E33").Value2 = "ok"
$ws.Range("F33").Value2 = 0.98317054260063
$ws.Range("G33").Value2 = 0
$ws.Range("H33").Value2 = 0.651006711409396
$ws.Range("I33").Value2 = 0
$ws.Range("J33").Value2 = 7.89
$ws.Range("K33").Value2 = 0.0196900269541779
$ws.Range("L33").Value2 = 0.0045754648011602
$ws.Range("M33").Value2 = 0.0353192169859886
$ws.Range("N33").Value2 = 0.249556742131532
$ws.Range("O33").Value2 = "RepSite"
$ws.Range("P33").Value2 = "Extremely likely increasing"
$ws.Range("Q33").Value2 = 1845196
$ws.Range("R33").Value2 = 5525095
$ws.Range("S33").Value2 = "Tararua District"
$ws.Range("T33").Value2 = "Manawatū"
$ws.Range("U33").Value2 = "Tiraumea"
$ws.Range("V33").Value2 = "Mana_7b"
$ws.Range("W33").Value2 = ""

# Row 34
$ws.Range("A34").Value2 = "Tiraumea u/s Manawatu Confluence"
$ws.Range("B34").Value2 = "SIN (Soluble Inorganic nitrogen)"
$ws.Range("C34").Value2 = 15
$ws.Range("D34").Value2 = $true
$ws.Range("E34").Value2 = "ok"
$ws.Range("F34").Value2 = 0.648819059783577
$ws.Range("G34").Value2 = 0
$ws.Range("H34").Value2 = 0.867549668874172
$ws.Range("I34").Value2 = 0
$ws.Range("J34").Value2 = 0.676
$ws.Range("K34").Value2 = -0.0015051510989011
$ws.Range("L34").Value2 = -0.0095371119579189
$ws.Range("M34").Value2 = 0.0043978321327234
$ws.Range("N34").Value2 = -0.22265548800313
$ws.Range("O34").Value2 = "RepSite"
$ws.Range("P34").Value2 = "As likely as not improving"
$ws.Range("Q34").Value2 = 1845196
$ws.Range("R34").Value2 = 5525095
$ws.Range("S34").Value2 = "Tararua District"
$ws.Range("T34").Value2 = "Manawatū"
$ws.Range("U34").Value2 = "Tiraumea"
$ws.Range("V34").Value2 = "Mana_7b"
$ws.Range("W34").Value2 = "g/m3"

# Row 35
$ws.Range("A35").Value2 = "Tiraumea u/s Manawatu Confluence"
$ws.Range("B35").Value2 = "Total Nitrogen"
$ws.Range("C35").Value2 = 15
$ws.Range("D35").Value2 = $true
$ws.Range("E35").Value2 = "ok"
$ws.Range("F35").Value2 = 0.997899099811342
$ws.Range("G35").Value2 = 0
$ws.Range("H35").Value2 = 0.622516556291391
$ws.Range("I35").Value2 = 0
$ws.Range("J35").Value2 = 0.86
$ws.Range("K35").Value2 = -0.0100343406593407
$ws.Range("L35").Value2 = -0.0166173794358508
$ws.Range("M35").Value2 = -0.0049818924068727
$ws.Range("N35").Value2 = -1.16678379759775
$ws.Range("O35").Value2 = "RepSite"
$ws.Range("P35").Value2 = "Virtually certain improving"
$ws.Range("Q35").Value2 = 1845196
$ws.Range("R35").Value2 = 5525095
$ws.Range("S35").Value2 = "Tararua District"
$ws.Range("T35").Value2 = "Manawatū"
$ws.Range("U35").Value2 = "Tiraumea"
$ws.Range("V35").Value2 = "Mana_7b"
$ws.Range("W35").Value2 = "g/m3"

# Row 36
$ws.Range("A36").Value2 = "Tiraumea u/s Manawatu Confluence"
$ws.Range("B36").Value2 = "Total Phosphorus"
$ws.Range("C36").Value2 = 15
$ws.Range("D36").Value2 = $true
$ws.Range("E36").Value2 = "ok"
$ws.Range("F36").Value2 = 0.989436092323271
$ws.Range("G36").Value2 = 0
$ws.Range("H36").Value2 = 0.384105960264901
$ws.Range("I36").Value2 = 0
$ws.Range("J36").Value2 = 0.024
$ws.Range("K36").Value2 = -0.0005001141031492
$ws.Range("L36").Value2 = -0.001003434065934
$ws.Range("M36").Value2 = -0.0001427817844367
$ws.Range("N36").Value2 = -2.08380876312186
$ws.Range("O36").Value2 = "RepSite"
$ws.Range("P36").Value2 = "Extremely likely improving"
$ws.Range("Q36").Value2 = 1845196
$ws.Range("R36").Value2 = 5525095
$ws.Range("S36").Value2 = "Tararua District"
$ws.Range("T36").Value2 = "Manawatū"
$ws.Range("U36").Value2 = "Tiraumea"
$ws.Range("V36").Value2 = "Mana_7b"
$ws.Range("W36").Value2 = "g/m3"

# Row 37
$ws.Range("A37").Value2 = "Tiraumea u/s Manawatu Confluence"
$ws.Range("B37").Value2 = "Turbidity"
$ws.Range("C37").Value2 = 15
$ws.Range("D37").Value2 = $true
$ws.Range("E37").Value2 = "ok"
$ws.Range("F37").Value2 = 0.998884336917713
$ws.Range("G37").Value2 = 0
$ws.Range("H37").Value2 = 0.940397350993378
$ws.Range("I37").Value2 = 0
$ws.Range("J37").Value2 = 4.37
$ws.Range("K37").Value2 = -0.169357178279197
$ws.Range("L37").Value2 = -0.288566836742088
$ws.Range("M37").Value2 = -0.0714397281370795
$ws.Range("N37").Value2 = -3.87545030387178
$ws.Range("O37").Value2 = "RepSite"
$ws.Range("P37").Value2 = "Virtually certain improving"
$ws.Range("Q37").Value2 = 1845196
$ws.Range("R37").Value2 = 5525095
$ws.Range("S37").Value2 = "Tararua District"
$ws.Range("T37").Value2 = "Manawatū"
$ws.Range("U37").Value2 = "Tiraumea"
$ws.Range("V37").Value2 = "Mana_7b"
$ws.Range("W37").Value2 = "NTU/FNU"

